$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy style from U1 then set values for V1:AA1
$ws.Range("U1").Copy()
$ws.Range("V1:AA1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(1, 22).Value = "Posesión Local ().2"
$ws.Cells.Item(1, 23).Value = "Posesión Visita ().2"
$ws.Cells.Item(1, 24).Value = "Posesión Local (%)"
$ws.Cells.Item(1, 25).Value = "Posesión Visita (%)"
$ws.Cells.Item(1, 26).Value = "fuente_tiempos"
$ws.Cells.Item(1, 27).Value = "estado_datos"


# Auto-generated corrections for rows 130-153 (M,N,O,P fixes)
$ws.Cells.Item(130, 14).Value = 1
$ws.Cells.Item(130, 16).Value = 2
$ws.Cells.Item(131, 13).Value = 1
$ws.Cells.Item(131, 15).Value = 0
$ws.Cells.Item(133, 14).Value = 1
$ws.Cells.Item(133, 16).Value = 0
$ws.Cells.Item(134, 13).Value = 1
$ws.Cells.Item(134, 15).Value = 1
$ws.Cells.Item(135, 13).Value = 1
$ws.Cells.Item(135, 14).Value = 2
$ws.Cells.Item(135, 15).Value = 0
$ws.Cells.Item(135, 16).Value = 1
$ws.Cells.Item(136, 13).Value = 2
$ws.Cells.Item(136, 15).Value = 2
$ws.Cells.Item(138, 13).Value = 1
$ws.Cells.Item(138, 15).Value = 1
$ws.Cells.Item(139, 13).Value = 1
$ws.Cells.Item(139, 15).Value = 1
$ws.Cells.Item(140, 13).Value = 2
$ws.Cells.Item(140, 14).Value = 1
$ws.Cells.Item(140, 15).Value = 1
$ws.Cells.Item(140, 16).Value = 0
$ws.Cells.Item(141, 14).Value = 1
$ws.Cells.Item(141, 16).Value = 1
$ws.Cells.Item(142, 13).Value = 1
$ws.Cells.Item(142, 15).Value = 1
$ws.Cells.Item(143, 13).Value = 2
$ws.Cells.Item(143, 14).Value = 1
$ws.Cells.Item(143, 15).Value = 1
$ws.Cells.Item(143, 16).Value = 1
$ws.Cells.Item(144, 14).Value = 1
$ws.Cells.Item(144, 16).Value = 1
$ws.Cells.Item(146, 13).Value = 3
$ws.Cells.Item(146, 14).Value = 1
$ws.Cells.Item(146, 15).Value = 1
$ws.Cells.Item(146, 16).Value = 3
$ws.Cells.Item(147, 13).Value = 1
$ws.Cells.Item(147, 14).Value = 1
$ws.Cells.Item(147, 15).Value = 1
$ws.Cells.Item(147, 16).Value = 0
$ws.Cells.Item(148, 13).Value = 2
$ws.Cells.Item(148, 14).Value = 1
$ws.Cells.Item(148, 15).Value = 0
$ws.Cells.Item(148, 16).Value = 0
$ws.Cells.Item(149, 14).Value = 4
$ws.Cells.Item(149, 16).Value = 0
$ws.Cells.Item(150, 13).Value = 1
$ws.Cells.Item(150, 15).Value = 1
$ws.Cells.Item(151, 13).Value = 1
$ws.Cells.Item(151, 14).Value = 1
$ws.Cells.Item(151, 15).Value = 1
$ws.Cells.Item(151, 16).Value = 1
$ws.Cells.Item(153, 14).Value = 1
$ws.Cells.Item(153, 16).Value = 3

# New rows 154-161 (full data including Q,R possession)
$ws.Range("A154:A163").NumberFormat = "@"
$ws.Cells.Item(154, 1).Value = "2025-08-08"
$ws.Cells.Item(154, 2).Value = "Dalian Zhixing"
$ws.Cells.Item(154, 3).Value = "Chengdu Better City"
$ws.Cells.Item(154, 4).Value = 0
$ws.Cells.Item(154, 5).Value = 2
$ws.Cells.Item(154, 6).Value = 1341051
$ws.Cells.Item(154, 7).Value = 4
$ws.Cells.Item(154, 8).Value = 2
$ws.Cells.Item(154, 9).Value = 3
$ws.Cells.Item(154, 10).Value = 6
$ws.Cells.Item(154, 11).Value = 0
$ws.Cells.Item(154, 12).Value = 1
$ws.Cells.Item(154, 13).Value = 0
$ws.Cells.Item(154, 14).Value = 0
$ws.Cells.Item(154, 15).Value = 0
$ws.Cells.Item(154, 16).Value = 2
$ws.Cells.Item(154, 17).Value = 53
$ws.Cells.Item(154, 18).Value = 47
$ws.Cells.Item(154, 19).Value = "V"
$ws.Cells.Item(155, 1).Value = "2025-08-08"
$ws.Cells.Item(155, 2).Value = "Henan Jianye"
$ws.Cells.Item(155, 3).Value = "Sichuan Jiuniu"
$ws.Cells.Item(155, 4).Value = 4
$ws.Cells.Item(155, 5).Value = 1
$ws.Cells.Item(155, 6).Value = 1341052
$ws.Cells.Item(155, 7).Value = 1
$ws.Cells.Item(155, 8).Value = 3
$ws.Cells.Item(155, 9).Value = 0
$ws.Cells.Item(155, 10).Value = 3
$ws.Cells.Item(155, 11).Value = 0
$ws.Cells.Item(155, 12).Value = 0
$ws.Cells.Item(155, 13).Value = 1
$ws.Cells.Item(155, 14).Value = 0
$ws.Cells.Item(155, 15).Value = 3
$ws.Cells.Item(155, 16).Value = 1
$ws.Cells.Item(155, 17).Value = 57
$ws.Cells.Item(155, 18).Value = 43
$ws.Cells.Item(155, 19).Value = "L"
$ws.Cells.Item(156, 1).Value = "2025-08-09"
$ws.Cells.Item(156, 2).Value = "Shandong Luneng"
$ws.Cells.Item(156, 3).Value = "Changchun Yatai"
$ws.Cells.Item(156, 4).Value = 2
$ws.Cells.Item(156, 5).Value = 1
$ws.Cells.Item(156, 6).Value = 1341053
$ws.Cells.Item(156, 7).Value = 6
$ws.Cells.Item(156, 8).Value = 7
$ws.Cells.Item(156, 9).Value = 3
$ws.Cells.Item(156, 10).Value = 3
$ws.Cells.Item(156, 11).Value = 1
$ws.Cells.Item(156, 12).Value = 0
$ws.Cells.Item(156, 13).Value = 1
$ws.Cells.Item(156, 14).Value = 0
$ws.Cells.Item(156, 15).Value = 1
$ws.Cells.Item(156, 16).Value = 1
$ws.Cells.Item(156, 17).Value = 60
$ws.Cells.Item(156, 18).Value = 40
$ws.Cells.Item(156, 19).Value = "L"
$ws.Cells.Item(157, 1).Value = "2025-08-09"
$ws.Cells.Item(157, 2).Value = "Shanghai Shenhua"
$ws.Cells.Item(157, 3).Value = "SHANGHAI SIPG"
$ws.Cells.Item(157, 4).Value = 1
$ws.Cells.Item(157, 5).Value = 2
$ws.Cells.Item(157, 6).Value = 1341054
$ws.Cells.Item(157, 7).Value = 9
$ws.Cells.Item(157, 8).Value = 4
$ws.Cells.Item(157, 9).Value = 4
$ws.Cells.Item(157, 10).Value = 6
$ws.Cells.Item(157, 11).Value = 0
$ws.Cells.Item(157, 12).Value = 0
$ws.Cells.Item(157, 13).Value = 0
$ws.Cells.Item(157, 14).Value = 1
$ws.Cells.Item(157, 15).Value = 1
$ws.Cells.Item(157, 16).Value = 1
$ws.Cells.Item(157, 17).Value = 54
$ws.Cells.Item(157, 18).Value = 46
$ws.Cells.Item(157, 19).Value = "V"
$ws.Cells.Item(158, 1).Value = "2025-08-09"
$ws.Cells.Item(158, 2).Value = "Qingdao Jonoon"
$ws.Cells.Item(158, 3).Value = "Yunnan Yukun"
$ws.Cells.Item(158, 4).Value = 5
$ws.Cells.Item(158, 5).Value = 1
$ws.Cells.Item(158, 6).Value = 1341055
$ws.Cells.Item(158, 7).Value = 8
$ws.Cells.Item(158, 8).Value = 1
$ws.Cells.Item(158, 9).Value = 2
$ws.Cells.Item(158, 10).Value = 2
$ws.Cells.Item(158, 11).Value = 0
$ws.Cells.Item(158, 12).Value = 1
$ws.Cells.Item(158, 13).Value = 2
$ws.Cells.Item(158, 14).Value = 1
$ws.Cells.Item(158, 15).Value = 3
$ws.Cells.Item(158, 16).Value = 0
$ws.Cells.Item(158, 17).Value = 60
$ws.Cells.Item(158, 18).Value = 40
$ws.Cells.Item(158, 19).Value = "L"
$ws.Cells.Item(159, 1).Value = "2025-08-10"
$ws.Cells.Item(159, 2).Value = "Tianjin Teda"
$ws.Cells.Item(159, 3).Value = "Qingdao Youth Island"
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 0
$ws.Cells.Item(159, 6).Value = 1341056
$ws.Cells.Item(159, 7).Value = 5
$ws.Cells.Item(159, 8).Value = 3
$ws.Cells.Item(159, 9).Value = 1
$ws.Cells.Item(159, 10).Value = 2
$ws.Cells.Item(159, 11).Value = 0
$ws.Cells.Item(159, 12).Value = 0
$ws.Cells.Item(159, 13).Value = 0
$ws.Cells.Item(159, 14).Value = 0
$ws.Cells.Item(159, 15).Value = 0
$ws.Cells.Item(159, 16).Value = 0
$ws.Cells.Item(159, 17).Value = 48
$ws.Cells.Item(159, 18).Value = 52
$ws.Cells.Item(159, 19).Value = "E"
$ws.Cells.Item(160, 1).Value = "2025-08-10"
$ws.Cells.Item(160, 2).Value = "Hangzhou Greentown"
$ws.Cells.Item(160, 3).Value = "Beijing Guoan"
$ws.Cells.Item(160, 4).Value = 3
$ws.Cells.Item(160, 5).Value = 4
$ws.Cells.Item(160, 6).Value = 1341057
$ws.Cells.Item(160, 7).Value = 5
$ws.Cells.Item(160, 8).Value = 6
$ws.Cells.Item(160, 9).Value = 1
$ws.Cells.Item(160, 10).Value = 2
$ws.Cells.Item(160, 11).Value = 1
$ws.Cells.Item(160, 12).Value = 0
$ws.Cells.Item(160, 13).Value = 2
$ws.Cells.Item(160, 14).Value = 1
$ws.Cells.Item(160, 15).Value = 1
$ws.Cells.Item(160, 16).Value = 3
$ws.Cells.Item(160, 17).Value = 35
$ws.Cells.Item(160, 18).Value = 65
$ws.Cells.Item(160, 19).Value = "V"
$ws.Cells.Item(161, 1).Value = "2025-08-10"
$ws.Cells.Item(161, 2).Value = "Wuhan Three Towns"
$ws.Cells.Item(161, 3).Value = "Meizhou Kejia"
$ws.Cells.Item(161, 4).Value = 1
$ws.Cells.Item(161, 5).Value = 2
$ws.Cells.Item(161, 6).Value = 1341058
$ws.Cells.Item(161, 7).Value = 3
$ws.Cells.Item(161, 8).Value = 3
$ws.Cells.Item(161, 9).Value = 1
$ws.Cells.Item(161, 10).Value = 3
$ws.Cells.Item(161, 11).Value = 0
$ws.Cells.Item(161, 12).Value = 0
$ws.Cells.Item(161, 13).Value = 0
$ws.Cells.Item(161, 14).Value = 1
$ws.Cells.Item(161, 15).Value = 1
$ws.Cells.Item(161, 16).Value = 1
$ws.Cells.Item(161, 17).Value = 68
$ws.Cells.Item(161, 18).Value = 32
$ws.Cells.Item(161, 19).Value = "V"

# New rows 162-163 (special: Q,R blank; X,Y,Z,AA populated)
$ws.Cells.Item(162, 1).Value = "2025-08-15"
$ws.Cells.Item(162, 2).Value = "SHANGHAI SIPG"
$ws.Cells.Item(162, 3).Value = "Henan Jianye"
$ws.Cells.Item(162, 4).Value = 4
$ws.Cells.Item(162, 5).Value = 1
$ws.Cells.Item(162, 6).Value = 1341059
$ws.Cells.Item(162, 7).Value = 4
$ws.Cells.Item(162, 8).Value = 4
$ws.Cells.Item(162, 9).Value = 3
$ws.Cells.Item(162, 10).Value = 3
$ws.Cells.Item(162, 11).Value = 0
$ws.Cells.Item(162, 12).Value = 0
$ws.Cells.Item(162, 13).Value = 3
$ws.Cells.Item(162, 14).Value = 0
$ws.Cells.Item(162, 15).Value = 1
$ws.Cells.Item(162, 16).Value = 1
$ws.Cells.Item(162, 19).Value = "L"
$ws.Cells.Item(162, 24).Value = 49
$ws.Cells.Item(162, 25).Value = 51
$ws.Cells.Item(162, 26).Value = "score"
$ws.Cells.Item(162, 27).Value = "OK"
$ws.Cells.Item(163, 1).Value = "2025-08-15"
$ws.Cells.Item(163, 2).Value = "Yunnan Yukun"
$ws.Cells.Item(163, 3).Value = "Wuhan Three Towns"
$ws.Cells.Item(163, 4).Value = 2
$ws.Cells.Item(163, 5).Value = 1
$ws.Cells.Item(163, 6).Value = 1341060
$ws.Cells.Item(163, 7).Value = 7
$ws.Cells.Item(163, 8).Value = 2
$ws.Cells.Item(163, 9).Value = 1
$ws.Cells.Item(163, 10).Value = 3
$ws.Cells.Item(163, 11).Value = 0
$ws.Cells.Item(163, 12).Value = 0
$ws.Cells.Item(163, 13).Value = 1
$ws.Cells.Item(163, 14).Value = 0
$ws.Cells.Item(163, 15).Value = 1
$ws.Cells.Item(163, 16).Value = 1
$ws.Cells.Item(163, 19).Value = "L"
$ws.Cells.Item(163, 24).Value = 50
$ws.Cells.Item(163, 25).Value = 50
$ws.Cells.Item(163, 26).Value = "score"
$ws.Cells.Item(163, 27).Value = "OK"

# Reset number format artifact on date cells so they keep default style (no explicit numFmt)
$ws.Range("B2").Copy()
$ws.Range("A154:A163").PasteSpecial(-4122)
$excel.CutCopyMode = $false


Write-Output "Edit complete"